$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.590.86"
$ws.Range("E2").Value = "  +2.41%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.648.78"
$ws.Range("E3").Value = "  +1.94%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.32"
$ws.Range("E5").Value = "  +1.42%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.32"
$ws.Range("E6").Value = "  +3.13%  "

# Row 8
$ws.Range("E8").Value = "  +0.12%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.647.76"
$ws.Range("E9").Value = "  +2.02%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.136"
$ws.Range("E10").Value = "  +9.34%  "

# Row 11
$ws.Range("E11").Value = "  -0.38%  "

# Row 12
$ws.Range("E12").Value = "  +1.29%  "

# Row 13
$ws.Range("E13").Value = "  +1.39%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.86"
$ws.Range("E14").Value = "  +2.33%  "

# Row 15
$ws.Range("E15").Value = "  +5.42%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.131.13"
$ws.Range("E16").Value = "  +1.94%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.477.24"
$ws.Range("E17").Value = "  +2.45%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.649.81"
$ws.Range("E18").Value = "  +2.00%  "

# Row 19
$ws.Range("E19").Value = "  +3.63%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "368.34"
$ws.Range("E20").Value = "  +1.31%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.43"
$ws.Range("E21").Value = "  +1.65%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.27"
$ws.Range("E22").Value = "  -0.20%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.86"
$ws.Range("E23").Value = "  +0.67%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.12"
$ws.Range("E24").Value = "  +3.77%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.95"
$ws.Range("E25").Value = "  +0.08%  "

# Row 26
$ws.Range("E26").Value = "  +0.02%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.97"
$ws.Range("E27").Value = "  +0.68%  "

# Row 28
$ws.Range("E28").Value = "  +7.08%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.770.66"
$ws.Range("E29").Value = "  +1.53%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.17%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "573.75"

# Row 32
$ws.Range("E32").Value = "  +4.69%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.00"
$ws.Range("E33").Value = "  +5.01%  "

# Row 34
$ws.Range("E34").Value = "  +2.66%  "

# Row 35
$ws.Range("E35").Value = "  +4.14%  "

# Row 36
$ws.Range("E36").Value = "  +0.10%  "

# Row 37
$ws.Range("E37").Value = "  +3.81%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "159.05"
$ws.Range("E38").Value = "  +1.94%  "

# Row 39
$ws.Range("E39").Value = "  +4.24%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.24"
$ws.Range("E40").Value = "  +1.55%  "

# Row 41
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.368"
$ws.Range("E41").Value = "  +0.84%  "

# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.39"
$ws.Range("E42").Value = "  +3.99%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.64"
$ws.Range("E43").Value = "  +5.11%  "

# Row 44
$ws.Range("E44").Value = "  +3.89%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₆0319"
$ws.Range("E45").Value = "  +12.52%  "

# Row 46
$ws.Range("E46").Value = "  +0.07%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.59"
$ws.Range("E47").Value = "  -0.41%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "157.20"
$ws.Range("E48").Value = "  +3.20%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.74"
$ws.Range("E49").Value = "  +0.93%  "

# Row 50
$ws.Range("E50").Value = "  +2.45%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.91"
$ws.Range("E51").Value = "  +2.66%  "
